$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title ---
Replace-Text "Unraveling the Enigma of Dark Energy" "Mathematics: The Universal Language of Science"

# --- Author ---
Replace-Text "Emma Jones" "John Roberts"

# --- Email address (emma.jones@hypothetical.edu -> johnroberts@eduworld.org) ---
Replace-Text "jones@hypothetical" ""
Replace-Text "edu" "org"
Replace-Text "emma" "johnroberts@eduworld"

# --- Body paragraph (dark energy -> mathematics) ---
Replace-Text "In the vast expanse of the cosmos, a perplexing enigma lingers - dark energy" "Mathematics, the cornerstone of scientific inquiry and technological advancement, has captivated human intellects for millennia"

Replace-Text " Its existence has been inferred from observations of the universe's expansion, an expansion that continues to accelerate" " It serves as a universal language, transcending linguistic and cultural boundaries, enabling us to comprehend and describe the intricate patterns inherent in the world around us"

Replace-Text " This mysterious force is believed to dominate the energy content of the universe, yet its nature remains veiled in obscurity" " Mathematics is the foundation upon which all scientific disciplines rest, providing a rigorous framework for analyzing, interpreting, and predicting natural phenomena"

Replace-Text " Our understanding of dark energy holds profound implications for our comprehension of the fundamental laws that govern the cosmos. Could it be a modification of gravity or an exotic field permeating space? Delving into the enigma of dark energy promises to unlock transformative insights into the workings of our universe" " Its influence is pervasive, touching every aspect of our lives, from the intricate workings of the atom to the vastness of the cosmos"

Replace-Text "Exploration of the accelerating expansion of the universe unveiled an enigma that continues to challenge our comprehension of physics: dark energy" "Through the lens of mathematics, we unravel the secrets of nature, unraveling the mysteries that have puzzled humanity for ages"

Replace-Text " This elusive entity is thought to wield the power to counteract the gravitational pull of matter, orchestrating the universe's ever-increasing expansion" " From the elegant simplicity of Pythagoras' Theorem to the complex intricacies of Einstein's Theory of Relativity, mathematics provides a roadmap, guiding us toward a deeper understanding of the universe"

Replace-Text " However, the mechanisms by which it exerts its influence remain shrouded in mystery. While some posit that dark energy is a cosmological constant, an unchanging property of spacetime, others propose dynamic scenarios, such as evolving scalar fields or modified theories of gravity. Determining the nature of dark energy holds the key to unlocking the secrets of the universe's ultimate fate" " It empowers us to unravel the enigmas of the quantum realm, probe the depths of black holes, and glimpse the birth and death of stars"

Replace-Text "The investigation into dark energy's enigmatic nature has yielded a plethora of theories, each attempting to decipher this perplexing phenomenon" "Furthermore, mathematics has revolutionized the way we live and work"

Replace-Text " One compelling idea suggests that dark energy might be a byproduct of quantum vacuum energy, the energy inherent in the vacuum state of space" " Its applications permeate diverse fields, including engineering, medicine, finance, and computer science"

Replace-Text " Another possibility involves a scalar field known as the `"quintessence field,`" which may permeate the universe and engender an accelerating expansion" " It underpins the development of life-saving technologies, facilitates global communication, and drives the innovation that shapes our modern world"

Replace-Text " Additionally, modifications to the conventional theory of gravity, such as the incorporation of additional dimensions or modifications to Einstein's equations, have also been proposed" " The profound impact of mathematics on society cannot be overstated"

Replace-Text " The pursuit of unraveling dark energy's true identity remains an ongoing endeavor, with scientists tirelessly seeking to illuminate this enigmatic chapter of cosmology" " It is the language of progress, the key to unlocking the boundless potential of human ingenuity"

# --- "Summary" heading: strip the stale lastRenderedPageBreak by rewriting the run in place ---
Replace-Text "Summary" "Summary"

# --- Summary paragraph body ---
Replace-Text "Our exploration into the nature of dark energy revealed the profound enigma it presents to our understanding of the universe" "Mathematics, the universal language of science, plays a pivotal role in comprehending the natural world, advancing technology, and driving societal progress"

Replace-Text " The accelerating expansion of the cosmos and the implications it holds for the future of our universe remain key areas of investigation" " Its influence spans diverse disciplines, from physics to biology to economics, providing a robust framework for analyzing and understanding complex phenomena"

Replace-Text " While various theories attempt to elucidate the mechanisms behind dark energy, its true identity remains shrouded in mystery" " Mathematics empowers us to unravel the secrets of nature, unraveling the mysteries that have puzzled humanity for ages"

Replace-Text " Unraveling this enigma promises to profoundly expand our comprehension of the fundamental laws governing the cosmos, offering transformative insights into the forces that orchestrate the evolution and ultimate fate of the universe" " It underpins the development of cutting-edge technologies, facilitates global communication, and fuels innovation. With its profound impact on both the scientific and social landscape, mathematics remains an indispensable tool for shaping a better future"

# --- Trailing empty paragraph at the end of the document ---
$count = $d.Paragraphs.Count
$d.Paragraphs($count).Range.InsertParagraphAfter() | Out-Null

Write-Host "Edit complete"
